# Adding 4 search test cases
# (2 new rows -> 4 new unique shared strings: test name + description for each row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Test Cases" sheet

# Seed the two new rows by copying the formatting (styles) of the last
# existing row (39), then overwrite the values.
$ws.Range("A39:E39").Copy($ws.Range("A40:E40"))
$ws.Range("A39:E39").Copy($ws.Range("A41:E41"))

# Row 40: PublishedAPostCommentCountTest
$ws.Cells.Item(40, 1).Value = "PublishedAPostCommentCountTest"
$ws.Cells.Item(40, 2).Value = "TBD"
$ws.Cells.Item(40, 3).Value = "Verify that POST tab Comment count getting increased while adding comment for post from Record view page"
$ws.Cells.Item(40, 4).Value = "Y"
$ws.Cells.Item(40, 5).Value = "SKIP"

# Row 41: OtherProfileWatchlistTabTest
$ws.Cells.Item(41, 1).Value = "OtherProfileWatchlistTabTest"
$ws.Cells.Item(41, 2).Value = "TBD"
$ws.Cells.Item(41, 3).Value = "Verify that user is able to watch others watchlists"
$ws.Cells.Item(41, 4).Value = "Y"
$ws.Cells.Item(41, 5).Value = "PASS"

# Reflect the updated scroll position / selection recorded in the sheet view.
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D2:D41").Select()
